# [#134279773] Importação de usuários aceita informação de turma
#
# Adds a new "Turma" (class/cohort) column to the user-import template:
#   - J1 header "Turma"
#   - J2 "Turma A" for the first sample user
#   - J3 "Turma B" for the second sample user
#   - J4 left blank for the third sample user
# Also nudges the active selection to the new next-available cell (J4),
# matching the author's final cursor position in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Turma" column - header + two sample values (third row intentionally left blank)
$ws.Range("J1").Value = "Turma"
$ws.Range("J2").Value = "Turma A"
$ws.Range("J3").Value = "Turma B"

# Move the selection to reflect where the editor left off after adding the column
$ws.Range("J4").Select()

# Minor cosmetic tweak carried over from the original commit (tab-bar / scrollbar split)
$excel.ActiveWindow.TabRatio = 985
